$wb = $excel.ActiveWorkbook

# Add the new "Result1" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Result1"

# Header row
$ws.Cells.Item(1,1).Value2 = "Symbol"
$ws.Cells.Item(1,2).Value2 = "LTP"
$ws.Cells.Item(1,3).Value2 = "High52W"
$ws.Cells.Item(1,4).Value2 = "DiffFrom52WHigh"
$ws.Cells.Item(1,5).Value2 = "DiffFrom52WLow"

# Data rows: Symbol, LTP, High52W, DiffFrom52WHigh, Low52W, DiffFrom52WLow
$data = @(
    @("TMPV", 367.4, 810.0, 45.358024691358025, 337.7, 108.79478827361564),
    @("TRENT", 4276.0, 7493.05, 57.06621469228152, 3930.1, 108.80130276583293),
    @("TCS", 3210.0, 4322.95, 74.2548491192357, 2866.6, 111.97934835693853),
    @("ETERNAL", 278.05, 368.45, 75.4647849097571, 194.8, 142.73613963039014),
    @("MAXHEALTH", 1045.1, 1314.3, 79.51761393897893, 940.05, 111.17493750332429),
    @("HCLTECH", 1620.1, 2012.2, 80.5138654209323, 1302.75, 124.36000767606984),
    @("INDIGO", 5060.0, 6232.5, 81.18732450862414, 3945.0, 128.2636248415716),
    @("WIPRO", 264.0, 324.6, 81.33086876155267, 228.0, 115.78947368421053),
    @("INFY", 1617.9, 1982.8, 81.59673189429091, 1307.0, 123.78729915837796),
    @("ITC", 402.7, 491.0, 82.01629327902239, 390.15, 103.21671152120979),
    @("POWERGRID", 264.2, 322.0, 82.04968944099377, 247.3, 106.83380509502626),
    @("HINDUNILVR", 2317.0, 2750.0, 84.25454545454545, 2136.0, 108.47378277153558),
    @("ADANIENT", 2239.0, 2616.5, 85.57232944773553, 1967.52, 113.79808083272343),
    @("APOLLOHOSP", 7039.0, 8099.5, 86.90659917278843, 6001.0, 117.29711714714215),
    @("JIOFIN", 295.0, 338.6, 87.12344949793265, 198.65, 148.50239114019632),
    @("ONGC", 240.46, 273.5, 87.91956124314443, 205.0, 117.29756097560977),
    @("NTPC", 329.9, 371.45, 88.81410687844932, 292.8, 112.6707650273224),
    @("BAJFINANCE", 986.5, 1102.5, 89.47845804988661, 679.2, 145.2444051825677),
    @("ICICIBANK", 1343.0, 1500.0, 89.53333333333333, 1186.0, 113.23777403035413),
    @("TECHM", 1591.0, 1773.6, 89.70455570590889, 1209.4, 131.55283611708285),
    @("ULTRACEMCO", 11792.0, 13097.0, 90.03588608078186, 10047.85, 117.3584398652448),
    @("SUNPHARMA", 1720.0, 1910.0, 90.0523560209424, 1548.0, 111.11111111111111),
    @("CIPLA", 1508.1, 1673.0, 90.14345487148833, 1335.0, 112.96629213483146),
    @("DRREDDY", 1273.0, 1405.9, 90.5469805818337, 1020.0, 124.80392156862746),
    @("HDFCLIFE", 749.05, 820.75, 91.2640877246421, 584.3, 128.19613212390894),
    @("BEL", 399.2, 436.0, 91.55963302752292, 240.25, 166.16024973985432),
    @("ASIANPAINT", 2770.2, 2985.7, 92.78226211608667, 2124.75, 130.37769149311683),
    @("BAJAJFINSV", 2038.0, 2195.0, 92.84738041002278, 1560.0, 130.64102564102564),
    @("ADANIPORTS", 1469.7, 1549.0, 94.88056810845707, 1010.75, 145.40687608211724),
    @("GRASIM", 2833.0, 2977.8, 95.13734972127072, 2276.95, 124.42082610509675),
    @("JSWSTEEL", 1165.9, 1223.9, 95.2610507394395, 880.0, 132.48863636363637),
    @("KOTAKBANK", 2203.0, 2301.9, 95.70354924193057, 1723.75, 127.8027556200145),
    @("COALINDIA", 399.45, 417.25, 95.73397243858597, 349.25, 114.37365783822476),
    @("TATASTEEL", 179.93, 186.94, 96.25013373274848, 122.62, 146.73788941445116),
    @("BHARTIARTL", 2107.0, 2174.5, 96.8958381237066, 1559.5, 135.10740621994228),
    @("AXISBANK", 1268.7, 1304.0, 97.29294478527608, 933.5, 135.90787359400107),
    @("HDFCBANK", 992.9, 1020.5, 97.29544341009309, 812.15, 122.25574093455643),
    @("SBILIFE", 2031.5, 2086.6, 97.35934055401131, 1372.55, 148.00917999344287),
    @("M&M", 3716.0, 3795.0, 97.91831357048747, 2425.0, 153.23711340206185),
    @("NESTLEIND", 1287.0, 1311.6, 98.12442817932296, 1055.0, 121.9905213270142),
    @("SBIN", 982.4, 999.0, 98.33833833833833, 680.0, 144.47058823529412),
    @("BAJAJ-AUTO", 9350.0, 9490.0, 98.52476290832455, 7089.35, 131.88797280427684),
    @("LT", 4083.0, 4140.0, 98.62318840579711, 2965.3, 137.69264492631436),
    @("TATACONSUM", 1190.1, 1202.8, 98.94413036248753, 893.1, 133.25495465233453),
    @("EICHERMOT", 7314.0, 7374.5, 99.17960539697606, 4646.0, 157.4257425742574),
    @("RELIANCE", 1569.0, 1581.3, 99.22215898311516, 1114.85, 140.73642194017134),
    @("MARUTI", 16688.0, 16818.0, 99.2270186704721, 10750.1, 155.23576524869534),
    @("HINDALCO", 886.0, 890.0, 99.5505617977528, 546.45, 162.13743251898617),
    @("TITAN", 4048.0, 4061.9, 99.65779561289052, 2925.0, 138.39316239316238),
    @("SHRIRAMFIN", 996.8, 1000.2, 99.66006798640271, 493.35, 202.04722813418462)
)

$rowCount = $data.Count
$colCount = 6
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range($ws.Cells.Item($startRow,1), $ws.Cells.Item($endRow,$colCount))
$rng.Value2 = $arr
